# Fill in the last benchmark row (row 16, Alveo U50 / MonteCarloPiEstimator)
# and add decimal precision to the power columns (K: FPGA power in Ws,
# L: FPGA on-chip power in W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing benchmark data for row 16 ---
# H16 ("FPGA utilization") was a blank placeholder; set its real value.
$ws.Range("H16").Value = 0.1622

# L16 ("FPGA on-chip power") was missing entirely; set its real value.
# K16 (=L16*J16/1000) and D16 (=(G16/K16)-1, previously #DIV/0!) will
# recompute automatically from this.
$ws.Range("L16").Value = 11.687

# --- Add two decimal places to the power columns ---
# Column L ("FPGA on-chip power"): was "0 ""W""", now "0.00 ""W""".
$ws.Range("L1:L17").NumberFormat = "0.00\ ""W"""

# Column K ("FPGA power"): was "0 ""Ws""", now "0.00 ""Ws""".
$ws.Range("K1:K17").NumberFormat = "0.00\ ""Ws"""

# Leave the selection on the recomputed cell (matches the saved cursor position).
[void]$ws.Range("D16").Select()
